$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D15").Value = 44243.52663644004
$ws.Range("D16:D29").Value = 44243.50540237268
$ws.Range("D30:D43").Value = 44243.48416751157
